$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns at D:E, shifting existing D:K to F:M
$ws.Columns("D:E").Insert()

# Copy number formats from column F (first surviving original column) into D:E
# so the new cells inherit the correct per-row style (date vs number).
$ws.Range("F7:F102").Copy()
$ws.Range("D7:E102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Rows 37 and 79 are section headers (only column B used) and had no D:K
# cells at all before the insert; PasteSpecial created blank placeholders
# there because the source column F also spans those rows. Remove them so
# those rows stay header-only, matching the original layout.
$ws.Range("D37:E37").Clear()
$ws.Range("D79:E79").Clear()


$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43373
$ws.Range("D8").Value = 1417100
$ws.Range("E8").Value = 1542200
$ws.Range("D9").Value = "NA"
$ws.Range("E9").Value = "NA"
$ws.Range("D10").Value = "NA"
$ws.Range("E10").Value = "NA"
$ws.Range("D12").Value = "NA"
$ws.Range("E12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("E14").Value = 0
$ws.Range("D15").Value = 33400
$ws.Range("E15").Value = 31700
$ws.Range("D17").Value = 1287200
$ws.Range("E17").Value = 1335800
$ws.Range("D18").Value = 129900
$ws.Range("E18").Value = 206400
$ws.Range("D20").Value = 0
$ws.Range("E20").Value = 0
$ws.Range("D21").Value = 163300
$ws.Range("E21").Value = 238100
$ws.Range("D22").Value = 11000
$ws.Range("E22").Value = 10800
$ws.Range("D23").Value = 118900
$ws.Range("E23").Value = 195600
$ws.Range("D24").Value = 25700
$ws.Range("E24").Value = 44100
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = 93200
$ws.Range("E26").Value = 151500
$ws.Range("D27").Value = 91600
$ws.Range("E27").Value = 151500
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = "NA"
$ws.Range("E29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = 0
$ws.Range("E32").Value = 0
$ws.Range("D33").Value = 91600
$ws.Range("E33").Value = 151500
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = 91600
$ws.Range("E35").Value = 151500
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43373
$ws.Range("D41").Value = 1467100
$ws.Range("E41").Value = 2205300
$ws.Range("D42").Value = 0
$ws.Range("E42").Value = 0
$ws.Range("D43").Value = 0
$ws.Range("E43").Value = 0
$ws.Range("D44").Value = 0
$ws.Range("E44").Value = 0
$ws.Range("D45").Value = 0
$ws.Range("E45").Value = 0
$ws.Range("D46").Value = 0
$ws.Range("E46").Value = 0
$ws.Range("D47").Value = 6638500
$ws.Range("E47").Value = 6637000
$ws.Range("D48").Value = 457800
$ws.Range("E48").Value = 457500
$ws.Range("D49").Value = 1253500
$ws.Range("E49").Value = 1256700
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("D52").Value = 16600
$ws.Range("E52").Value = 22800
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 10630600
$ws.Range("E54").Value = 11380100
$ws.Range("D57").Value = 778700
$ws.Range("E57").Value = 826700
$ws.Range("D58").Value = 0
$ws.Range("E58").Value = 0
$ws.Range("D59").Value = 4837900
$ws.Range("E59").Value = 5581000
$ws.Range("D60").Value = 0
$ws.Range("E60").Value = 0
$ws.Range("D61").Value = 808300
$ws.Range("E61").Value = 821800
$ws.Range("D62").Value = 217100
$ws.Range("E62").Value = 219300
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("D66").Value = 6888800
$ws.Range("E66").Value = 7705900
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("D72").Value = 1644200
$ws.Range("E72").Value = 1600300
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("D76").Value = 3741900
$ws.Range("E76").Value = 3674200
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43373
$ws.Range("D81").Value = 91600
$ws.Range("E81").Value = 151500
$ws.Range("D83").Value = 33400
$ws.Range("E83").Value = 31700
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("D89").Value = 308300
$ws.Range("E89").Value = 230800
$ws.Range("D91").Value = -30900
$ws.Range("E91").Value = -31600
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("D94").Value = -205900
$ws.Range("E94").Value = -446800
$ws.Range("D96").Value = -46900
$ws.Range("E96").Value = -46900
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("D100").Value = -836700
$ws.Range("E100").Value = 1197000
$ws.Range("D101").Value = -3900
$ws.Range("E101").Value = -2200
$ws.Range("D102").Value = -738200
$ws.Range("E102").Value = 978800
